$p = $ppt.ActivePresentation

# Slide 13 ("...Slides.Item(13)") holds the bottom banner shape "Rectangle 5"
# (cNvPr id="6") that names the teacher / school in Arabic. It is the last
# shape in the slide's shape tree and is being removed entirely.
$s = $p.Slides.Item(13)
$shapes = $s.Shapes
for ($i = $shapes.Count; $i -ge 1; $i--) {
    $shp = $shapes.Item($i)
    if ($shp.Name -eq "Rectangle 5") {
        $shp.Delete()
    }
}
